# Add simulation results for rollout
$wb = $excel.ActiveWorkbook

# --- raw_scores sheet: add I (rollout) and J (rollout_smooth) columns ---
$wsRaw = $wb.Worksheets.Item("raw_scores")
$wsRaw.Range("I1").Value = "rollout"
$wsRaw.Range("J1").Value = "rollout_smooth"

$iVals = @(5380,5568,6884,6560,3288,1056,1340,3616,11464,6328,5348,1352,4484,1180,6596,4880,2280,1404,6792,2680,3476,2432,5988,2452,6276,7112,1264,4940,7072,1352,5312,1552,1172,7380,6356,5612,3124,1108,3640,4940,1308,6328,3232,1428,1612,2448,5668,5572,6696,3384,2836,6484,3108,3236,3200,3196,764,1460,5304,6756,11912,2988,3008,5396,4528,3100,3248,3256,2500,6464,7080,2344,6324,5504,1324,5356,3096,3160,3064,6136,3764,5176,3144,6884,10636,1016,7096,3240,2472,5120,3468,3520,7332,1472,3144,5660,6396,6928,6448,3520)
$jVals = @(3140,1388,5396,2916,3108,5572,3060,3108,1336,736,6424,5264,780,304,3332,988,5072,1428,3488,6028,5440,2672,2348,2388,1436,6204,2692,1372,3328,1468,4944,2460,5416,6284,4084,7188,3196,2684,5228,6864,4944,1544,5628,1452,792,6172,2936,2232,4576,1180,6536,1316,3196,3516,1720,6016,4556,1448,5172,936,7312,3116,2316,5380,3508,1680,7232,3116,5424,1520,7232,2372,1016,2388,2596,5420,4620,2808,2748,7248,4900,3160,6908,7048,3132,4916,2912,6056,4532,844,3116,5300,2492,5124,4584,5204,1524,6632,1724,1924)
$data = New-Object 'object[,]' 100,2
for ($row = 0; $row -lt 100; $row++) {
    $data[$row,0] = $iVals[$row]
    $data[$row,1] = $jVals[$row]
}
$wsRaw.Range("I2:J101").Value = $data

# Row 2 (first data row) carries the Menlo-font style present on A2:H2;
# copy that formatting onto the new I2:J2 cells so they match.
$wsRaw.Range("H2").Copy()
$wsRaw.Range("I2:J2").PasteSpecial(-4122)
$wsRaw.Range("A1").Select()

# --- hyper sheet: rename MCTS -> MCTS_default, add rollout_default / rollout_smooth columns + machine specs ---
$wsHyper = $wb.Worksheets.Item("hyper")
$wsHyper.Range("I1").Value = "MCTS_default"
$wsHyper.Range("J1").Value = "rollout_default"
$wsHyper.Range("K1").Value = "rollout_smooth"
$wsHyper.Range("J2").Value = 12948
$wsHyper.Range("K2").Value = 15214
$wsHyper.Range("J3").Value = 100
$wsHyper.Range("K3").Value = 100

$wsHyper.Range("A5").Value = "Processor"
$wsHyper.Range("B5").Value = "2.2 GHz Intel Core i7"
$wsHyper.Range("A6").Value = "Memory"
$wsHyper.Range("B6").Value = "16 GB 1600 MHz DDR3"
$wsHyper.Range("A7").Value = "Graphics"
$wsHyper.Range("B7").Value = "Intel Iris Pro 1536 MB"
$wsHyper.Range("A1").Select()

# --- statistics sheet: rename mcts header source, add rollout (J), rollout_smooth (K), rollout_depth5 (L, header only) ---
$wsStat = $wb.Worksheets.Item("statistics")
$wsStat.Range("J1").Value = "rollout"
$wsStat.Range("K1").Value = "rollout_smooth"
$wsStat.Range("L1").Value = "rollout_depth5"

$wsStat.Range("J2").Formula = "=MIN(raw_scores!I2:raw_scores!I101)"
$wsStat.Range("K2").Formula = "=MIN(raw_scores!J2:raw_scores!J101)"

$wsStat.Range("J3").Formula = "=MAX(raw_scores!I2:raw_scores!I101)"
$wsStat.Range("K3").Formula = "=MAX(raw_scores!J2:raw_scores!J101)"

$wsStat.Range("J4").Formula = "=AVERAGE(raw_scores!I2:raw_scores!I101)"
$wsStat.Range("K4").Formula = "=AVERAGE(raw_scores!J2:raw_scores!J101)"

$wsStat.Range("J5").Formula = "=MEDIAN(raw_scores!I2:raw_scores!I101)"
$wsStat.Range("K5").Formula = "=MEDIAN(raw_scores!J2:raw_scores!J101)"

$wsStat.Range("J6").Formula = "=STDEV(raw_scores!I2:raw_scores!I101)"
$wsStat.Range("K6").Formula = "=STDEV(raw_scores!J2:raw_scores!J101)"

$wsStat.Range("J7").Formula = "=hyper!J2/hyper!J3"
$wsStat.Range("K7").Formula = "=hyper!K2/hyper!K3"
$wsStat.Range("A1").Select()

$wsHyper.Activate()
$wsHyper.Range("A1").Select()
